# Updates the cached "today" date shown by the datetimeFigureOut footer
# field on the slide master and every slide layout, and registers a
# custom show ("Custom Show 1") that contains the presentation's single
# slide.

function Set-DatePlaceholderText($Shapes, $NewText) {
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $shp = $Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePlaceholder = $false
            try {
                # ppPlaceholderDate = 16
                $isDatePlaceholder = ($shp.PlaceholderFormat.Type -eq 16)
            } catch {
                $isDatePlaceholder = $false
            }

            if ($isDatePlaceholder) {
                $shp.TextFrame.TextRange.Text = $NewText
            }
        }
    }
}

$p = $ppt.ActivePresentation
$newDate = "2/13/2021"

# Slide master footer date placeholder.
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every slide layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Register a custom show containing the (only) slide of the deck.
$slideIds = @()
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slideIds += $p.Slides.Item($si).SlideID
}

$namedShows = $p.SlideShowSettings.NamedSlideShows
$namedShows.Add("Custom Show 1", $slideIds) | Out-Null
